# Build and QA Cleanup
#
# 1) Refresh the cached "datetimeFigureOut" date field text (footer date)
#    from 4/4/2025 to 8/13/2025 across the slide master, every slide
#    layout, and the notes master.
# 2) Clean up the "HRex Patient" labels (drop the stray leading "HRex"
#    run and rename the remaining label to "US Core Patient") wherever
#    they occur, including inside grouped shapes.

$newDate = "8/13/2025"

function Set-DatePlaceholderText($container) {
    for ($i = 1; $i -le $container.Shapes.Count; $i++) {
        $shape = $container.Shapes.Item($i)
        if ($shape.Name -like "Date Placeholder*") {
            if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
                $shape.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

$p = $ppt.ActivePresentation

# --- Slide master date placeholder ---
$master = $p.SlideMaster
Set-DatePlaceholderText $master

# --- Every slide layout's date placeholder ---
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Set-DatePlaceholderText $layout
}

# --- Notes master date placeholder ---
$notesMaster = $p.NotesMaster
Set-DatePlaceholderText $notesMaster

# --- Fix up "HRex Patient" -> "US Core Patient" labels ---
function Fix-HRexPatientLabel($shape) {
    if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
        $tr = $shape.TextFrame.TextRange
        if ($tr.Text -eq "HRex Patient") {
            # Keep the formatting of the " Patient" run (2nd run) and
            # expand it to the full new label...
            $tr.Characters(5, 8).Text = "US Core Patient"
            # ...then drop the stray leading "HRex" run entirely.
            $tr.Characters(1, 4).Text = ""
        }
    }
    if ($shape.Type -eq 6) {
        for ($j = 1; $j -le $shape.GroupItems.Count; $j++) {
            Fix-HRexPatientLabel $shape.GroupItems.Item($j)
        }
    }
}

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        Fix-HRexPatientLabel $slide.Shapes.Item($i)
    }
}
